# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" between "2021-Q4" and "总计", built as a
#   structural copy of "2021-Q4" (same header row/style) trimmed down to a
#   single data row for fund 005021.
# - Update the "总计" sheet: insert a new top data row for "2022-Q1"
#   (count=1, total=0) and renumber the existing index column.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to store a literal text value (keeps "005021" etc. from
    # being coerced into a number and losing the leading zero), then drop
    # back to the default "Normal" style so no stray number-format style
    # sticks to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws2021Q3 = $wb.Worksheets.Item(1)
$ws2021Q4 = $wb.Worksheets.Item(2)
$wsTotal  = $wb.Worksheets.Item(3)

# --- 1. Create the "2022-Q1" sheet right after "2021-Q4" ------------------
$ws2021Q4.Copy($null, $ws2021Q4)
$ws2022Q1 = $wb.Worksheets.Item(3)
$ws2022Q1.Name = "2022-Q1"

# The copied sheet has 4 data rows (rows 2-5); keep only the first one.
$ws2022Q1.Rows.Item(3).EntireRow.Delete()
$ws2022Q1.Rows.Item(3).EntireRow.Delete()
$ws2022Q1.Rows.Item(3).EntireRow.Delete()

Set-TextValue $ws2022Q1.Range("B2") "005021"
Set-TextValue $ws2022Q1.Range("C2") "渤海汇金量化汇盈灵活配置混合"
Set-TextValue $ws2022Q1.Range("D2") "0.02"
Set-TextValue $ws2022Q1.Range("E2") "92.66"
Set-TextValue $ws2022Q1.Range("F2") "1.15"
Set-TextValue $ws2022Q1.Range("G2") "0.0002"
$ws2022Q1.Range("H2").Value = 8

# --- 2. Update the "总计" sheet --------------------------------------------
$wsTotal = $wb.Worksheets.Item(4)

$wsTotal.Rows.Item(2).EntireRow.Insert()

# Inserting a row copies the formatting of the row above into the new row;
# re-paint it from the (now shifted-down) old first data row so the new
# row2 matches the plain data-row style instead.
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0

# Renumber the 0-based index column for the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# --- 3. Restore the originally active sheet/selection ---------------------
$ws2021Q3.Select()

Write-Host "done"
